# Commit message: "added raw accuracy and nan"
#
# The original sheet had columns: A=Attack, B=Epsilon, C..H = detector
# success rates (BAARD2, BAARD3, FS, LID, RC, magnet).
#
# This edit inserts a new column before the old column C, shifting the
# detector columns from C:H to D:I, and uses the freed-up column C to
# hold the raw "Accuracy after attack" values (with header in C3 and
# the same header-row styling as the existing B1/B2 blank header cells
# in rows 1-2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; this shifts the previous C:H columns
# (BAARD2, BAARD3, FS, LID, RC, magnet data) one column to the right,
# to D:I, and copies formatting from the neighboring column.
$ws.Columns("C").Insert()

# Rows 1-2 hold the two-row "group" header (blank cell over the
# Attack/Epsilon columns, then the detector names / "Success" row);
# that blank leading cell was only ever one cell wide and simply
# slides over to the newly inserted column, so the stray leftover in
# column B (duplicated by the column insert above) needs clearing.
$ws.Range("B1").Clear()
$ws.Range("B2").Clear()

# New header for the inserted column.
$ws.Range("C3").Value = "Accuracy after attack"

# Raw accuracy-after-attack values per row (rows 4-33).
$ws.Range("C4").Value = 0.1
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 99.40000000000001
$ws.Range("C10").Value = 98.90000000000001
$ws.Range("C11").Value = 98.3
$ws.Range("C12").Value = 96.59999999999999
$ws.Range("C13").Value = 94.8
$ws.Range("C14").Value = 91.3
$ws.Range("C15").Value = 82.5
$ws.Range("C16").Value = 46.7
$ws.Range("C17").Value = 17.4
$ws.Range("C18").Value = 12.9
$ws.Range("C19").Value = 9.800000000000001
$ws.Range("C20").Value = 6.600000000000001
$ws.Range("C21").Value = 3.1
$ws.Range("C22").Value = 0.6
$ws.Range("C23").Value = 0.2
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("C28").Value = 3.7
$ws.Range("C29").Value = 8.4
$ws.Range("C30").Value = 7.5
$ws.Range("C31").Value = 8.4
$ws.Range("C32").Value = 7.000000000000001
$ws.Range("C33").Value = 7.000000000000001
